# Update countries & provincias Spain
# Refresh the "last updated" timestamp and the per-country COVID figures.
# A handful of countries' case counts grew enough to overtake their
# neighbours in the (descending, by total cases) ordering, so several rows
# swap their country label as well as their numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner in row 1
$ws.Range("A1").Value = "Datos actualizados a 22 de Abril de 2020 a las 10:22"

# Helper: write a full A:H row, cell by cell (column 1 = A, ... column 8 = H)
function Set-RowValues($row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# Polonia keeps its position but gets refreshed counts
Set-RowValues 32 @("Polonia", 10034, 178, 1513, 8117, 160, 3, 404)

# Australia: only "Muertes hoy" (F) changes
$ws.Cells.Item(43, 6).Value = 47

# Eslovaquia overtakes Republica de Macedonia
Set-RowValues 78 @("Eslovaquia", 1244, 45, 284, 946, 7, 0, 14)
Set-RowValues 79 @("Republica de Macedonia", 1231, 0, 224, 952, 9, 0, 55)

# Letonia keeps its position but gets refreshed counts
Set-RowValues 91 @("Letonia", 761, 13, 133, 617, 5, 2, 11)

# Sri Lanka overtakes Guatemala, Montenegro and Mayotte
Set-RowValues 114 @("Sri Lanka", 321, 11, 104, 210, 2, 0, 7)
Set-RowValues 115 @("Guatemala", 316, 22, 24, 284, 3, 1, 8)
Set-RowValues 116 @("Montenegro", 314, 1, 101, 208, 7, 0, 5)
Set-RowValues 117 @("Mayotte", 311, 0, 117, 190, 4, 0, 4)

# Malaui overtakes Botsuana, Laos and Belice
Set-RowValues 180 @("Malaui", 23, 5, 3, 18, 1, 0, 2)
Set-RowValues 181 @("Botsuana", 20, 0, 0, 19, 0, 0, 1)
Set-RowValues 182 @("Laos", 19, 0, 4, 15, 0, 0, 0)
Set-RowValues 183 @("Belice", 18, 0, 2, 13, 1, 0, 2)
